# Weekly fruit/vegetable price update.
# Insert a new daily-price record as row 84 (pushing the existing rows 84-102
# down to 85-103) for "Berenjena" sourced from Provincia de Huasco.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 84; this shifts the prior rows 84-102
# down to 85-103 and enlarges the used range to A1:R103.
$ws.Rows.Item(84).EntireRow.Insert()

# Populate the newly inserted row 84 with the new record.
$ws.Range("A84").Value = 6
$ws.Range("B84").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C84").Value = "Metropolitana"
$ws.Range("D84").Value = 44476
$ws.Range("E84").Value = 13
$ws.Range("F84").Value = 100112001
$ws.Range("G84").Value = "Berenjena"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 220
$ws.Range("K84").Value = 7000
$ws.Range("L84").Value = 8000
$ws.Range("M84").Value = 7545
$ws.Range("N84").Value = "$/caja 60 unidades"
$ws.Range("O84").Value = "Provincia de Huasco"
$ws.Range("P84").Value = 126
$ws.Range("Q84").Value = 60
$ws.Range("R84").Value = "Hortaliza"
